# Edit script: applies the diff to before.xlsx
$wb = $excel.ActiveWorkbook

# Sheet "Resumen" - summary values
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("B2").Value = "Z3"
$wsResumen.Range("C2").Value = 644.2071305005946

# Sheet "Solucion" - the randomized/mutated pedido->salida assignment
$wsSolucion = $wb.Worksheets.Item("Solucion")
$solucionData = @(
    @(2, "Pedido_23", "S001"),
    @(3, "Pedido_42", "S025"),
    @(4, "Pedido_60", "S041"),
    @(5, "Pedido_5", "S045"),
    @(6, "Pedido_11", "S029"),
    @(7, "Pedido_15", "S005"),
    @(8, "Pedido_16", "S042"),
    @(9, "Pedido_56", "S002"),
    @(10, "Pedido_22", "S026"),
    @(11, "Pedido_44", "S006"),
    @(12, "Pedido_55", "S046"),
    @(13, "Pedido_4", "S030"),
    @(14, "Pedido_51", "S043"),
    @(15, "Pedido_40", "S003"),
    @(16, "Pedido_25", "S047"),
    @(17, "Pedido_1", "S007"),
    @(18, "Pedido_7", "S027"),
    @(19, "Pedido_3", "S004"),
    @(20, "Pedido_38", "S044"),
    @(21, "Pedido_49", "S031"),
    @(22, "Pedido_35", "S008"),
    @(23, "Pedido_33", "S048"),
    @(24, "Pedido_48", "S028"),
    @(25, "Pedido_31", "S049"),
    @(26, "Pedido_30", "S009"),
    @(27, "Pedido_59", "S053"),
    @(28, "Pedido_54", "S032"),
    @(29, "Pedido_2", "S013"),
    @(30, "Pedido_45", "S050"),
    @(31, "Pedido_29", "S033"),
    @(32, "Pedido_24", "S010"),
    @(33, "Pedido_58", "S054"),
    @(34, "Pedido_17", "S037"),
    @(35, "Pedido_36", "S014"),
    @(36, "Pedido_57", "S051"),
    @(37, "Pedido_37", "S034"),
    @(38, "Pedido_39", "S011"),
    @(39, "Pedido_34", "S055"),
    @(40, "Pedido_52", "S052"),
    @(41, "Pedido_6", "S015"),
    @(42, "Pedido_27", "S038"),
    @(43, "Pedido_47", "S056"),
    @(44, "Pedido_46", "S012"),
    @(45, "Pedido_10", "S035"),
    @(46, "Pedido_19", "S016"),
    @(47, "Pedido_50", "S057"),
    @(48, "Pedido_28", "S061"),
    @(49, "Pedido_26", "S017"),
    @(50, "Pedido_21", "S039"),
    @(51, "Pedido_12", "S036"),
    @(52, "Pedido_20", "S058"),
    @(53, "Pedido_41", "S021"),
    @(54, "Pedido_32", "S040"),
    @(55, "Pedido_13", "S062"),
    @(56, "Pedido_18", "S018"),
    @(57, "Pedido_9", "S059"),
    @(58, "Pedido_8", "S022"),
    @(59, "Pedido_53", "S063"),
    @(60, "Pedido_14", "S019"),
    @(61, "Pedido_43", "S060")
)

foreach ($row in $solucionData) {
    $r = $row[0]
    $pedido = $row[1]
    $salida = $row[2]
    $wsSolucion.Cells.Item($r, 1).Value = $pedido
    $wsSolucion.Cells.Item($r, 2).Value = $salida
}

# Sheet "Metricas" - per zone time metrics
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 643.1801275813602
$wsMetricas.Range("B3").Value = 508.8339820521136
$wsMetricas.Range("B4").Value = 644.2071305005946
